# Bump the "Förändrad" date column (C2:C171) from 45189 to 45190 (2023-09-20 -> 2023-09-21)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 171
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45190
